# Fix #4005 Model ISA-TAB-MODEL.xlsx import error
# Adds a machine-readable "name" per attribute row (short identifier) while
# moving the previous descriptive "name" text into a new "label" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Row 2..131: new-name (becomes column A) -> old-name (becomes column B / "label")
$rowData = @{
    2 = @('inv_Identifier','Investigation Identifier')
    3 = @('inv_Title','Investigation Title')
    4 = @('inv_Description','Investigation Description')
    5 = @('inv_Submission_Date','Investigation Submission Date')
    6 = @('inv_Public_Release_Date','Investigation Public Release Date')
    7 = @('Term_src_name','Term source name')
    8 = @('Term_src_File','Term Source File')
    9 = @('Term_src_Version','Term Source Version')
    10 = @('Term_src_Description','Term Source Description')
    11 = @('inv_PubMed_ID','Investigation PubMed ID')
    12 = @('inv_pub_DOI','Investigation Publication DOI')
    13 = @('inv_pub_Author_List','Investigation Publication Author List')
    14 = @('inv_pub_Title','Investigation Publication Title')
    15 = @('inv_pub_Status','Investigation Publication Status')
    16 = @('inv_pub_Status_Term_acc_num','Investigation Publication Status Term Accession Number')
    17 = @('inv_pub_Status_Term_src_REF','Investigation Publication Status Term Source REF')
    18 = @('inv_Person_Last_Name','Investigation Person Last Name')
    19 = @('inv_Person_First_Name','Investigation Person First Name')
    20 = @('inv_Person_Mid_Initials','Investigation Person Mid Initials')
    21 = @('inv_Person_Email','Investigation Person Email')
    22 = @('inv_Person_Phone','Investigation Person Phone')
    23 = @('inv_Person_Fax','Investigation Person Fax')
    24 = @('inv_Person_Address','Investigation Person Address')
    25 = @('inv_Person_Affiliation','Investigation Person Affiliation')
    26 = @('inv_Person_Roles_','Investigation Person Roles ')
    27 = @('inv_Person_Roles_Term_acc_num','Investigation Person Roles Term Accession Number')
    28 = @('inv_Person_Roles_Term_src_REF','Investigation Person Roles Term Source REF')
    29 = @('inv_Identifier','Investigation Identifier')
    30 = @('st_Identifier','Study Identifier')
    31 = @('st_Title','Study Title')
    32 = @('st_Description','Study Description')
    33 = @('st_Submission_Date','Study Submission Date')
    34 = @('st_Public_Release_Date','Study Public Release Date')
    35 = @('st_File_Name','Study File Name')
    36 = @('st_Design_tp','Study Design Type')
    37 = @('st_Design_tp_Term_acc_num','Study Design Type Term Accession Number')
    38 = @('st_Design_tp_Term_src_REF','Study Design Type Term Source REF')
    39 = @('st_PubMed_ID','Study PubMed ID')
    40 = @('st_pub_DOI','Study Publication DOI')
    41 = @('st_pub_Author_List','Study Publication Author List')
    42 = @('st_pub_Title','Study Publication Title')
    43 = @('st_pub_Status','Study Publication Status')
    44 = @('st_pub_Status_Term_acc_num','Study Publication Status Term Accession Number')
    45 = @('st_pub_Status_Term_src_REF','Study Publication Status Term Source REF')
    46 = @('st_Factor_Name','Study Factor Name')
    47 = @('st_Factor_tp','Study Factor Type')
    48 = @('st_Factor_tp_Term_acc_num','Study Factor Type Term Accession Number')
    49 = @('st_Factor_tp_Term_src_REF','Study Factor Type Term Source REF')
    50 = @('st_Assay_meas_tp','Study Assay Measurement Type')
    51 = @('st_Assay_meas_tp_Term_acc_num','Study Assay Measurement Type Term Accession Number')
    52 = @('st_Assay_meas_tp_Term_src_REF','Study Assay Measurement Type Term Source REF')
    53 = @('st_Assay_tech_tp','Study Assay Technology Type')
    54 = @('st_Assay_tech_tp_Term_acc_num','Study Assay Technology Type Term Accession Number')
    55 = @('st_Assay_tech_tp_Term_src_REF','Study Assay Technology Type Term Source REF')
    56 = @('st_Assay_tech_Platform','Study Assay Technology Platform')
    57 = @('st_Assay_File_Name','Study Assay File Name')
    58 = @('st_prot_Name','Study Protocol Name')
    59 = @('st_prot_tp','Study Protocol Type')
    60 = @('st_prot_tp_Term_acc_num','Study Protocol Type Term Accession Number')
    61 = @('st_prot_tp_Term_src_REF','Study Protocol Type Term Source REF')
    62 = @('st_prot_Description','Study Protocol Description')
    63 = @('st_prot_URI','Study Protocol URI')
    64 = @('st_prot_Version','Study Protocol Version')
    65 = @('st_prot_params_Name','Study Protocol Parameters Name')
    66 = @('st_prot_params_Term_acc_num','Study Protocol Parameters Term Accession Number')
    67 = @('st_prot_params_Term_src_REF','Study Protocol Parameters Term Source REF')
    68 = @('st_prot_comps_Name','Study Protocol Components Name')
    69 = @('st_prot_comps_tp','Study Protocol Components Type')
    70 = @('st_prot_comps_tp_Term_acc_num','Study Protocol Components Type Term Accession Number')
    71 = @('st_prot_comps_tp_Term_src_REF','Study Protocol Components Type Term Source REF')
    72 = @('st_Person_Last_Name','Study Person Last Name')
    73 = @('st_Person_First_Name','Study Person First Name')
    74 = @('st_Person_Mid_Initials','Study Person Mid Initials')
    75 = @('st_Person_Email','Study Person Email')
    76 = @('st_Person_Phone','Study Person Phone')
    77 = @('st_Person_Fax','Study Person Fax')
    78 = @('st_Person_Address','Study Person Address')
    79 = @('st_Person_Affiliation','Study Person Affiliation')
    80 = @('st_Person_Roles_','Study Person Roles ')
    81 = @('st_Person_Roles_Term_acc_num','Study Person Roles Term Accession Number')
    82 = @('st_Person_Roles_Term_src_REF','Study Person Roles Term Source REF')
    83 = @('st_Identifier','Study Identifier')
    84 = @('Source_Name','Source Name')
    85 = @('Sample_Name','Sample Name')
    86 = @('Material_tp','Material Type')
    87 = @('Characteristics','Characteristics')
    88 = @('Provider','Provider')
    89 = @('Protocol_REF','Protocol REF')
    90 = @('Term_acc_num','Term Accession Number')
    91 = @('Term_src_REF','Term Source REF')
    92 = @('Factor_Value','Factor Value')
    93 = @('Comment','Comment')
    94 = @('Sample_Name_','Sample Name ')
    95 = @('Extract_Name_','Extract Name ')
    96 = @('Labeled_Extract_Name_','Labeled Extract Name ')
    97 = @('Image_File','Image File')
    98 = @('Raw_Data_File','Raw Data File')
    99 = @('Data_Transformation_Name','Data Transformation Name')
    100 = @('Normalization_Name','Normalization Name')
    101 = @('Derived_Data_File','Derived Data File')
    102 = @('Material_tp','Material Type')
    103 = @('Characteristics','Characteristics')
    104 = @('Label','Label')
    105 = @('Protocol_REF','Protocol REF')
    106 = @('Term_acc_num','Term Accession Number')
    107 = @('Term_src_REF','Term Source REF')
    108 = @('Comment','Comment')
    109 = @('Hybridization_Assay_Name','Hybridization Assay Name')
    110 = @('Scan_Name','Scan Name')
    111 = @('Array_Data_File','Array Data File')
    112 = @('Derived_Array_Data_File','Derived Array Data File')
    113 = @('Array_Data_Matrix_File','Array Data Matrix File')
    114 = @('Derived_Array_Data_Matrix_File','Derived Array Data Matrix File')
    115 = @('Array_Design_File','Array Design File')
    116 = @('Array_Design_REF','Array Design REF')
    117 = @('Gel_Electrophoresis_Assay_Name','Gel Electrophoresis Assay Name')
    118 = @('First_Dimension_','First Dimension ')
    119 = @('Second_Dimension','Second Dimension')
    120 = @('Scan_Name','Scan Name')
    121 = @('Spot_Picking_File','Spot Picking File')
    122 = @('MS_Assay_Name','MS Assay Name')
    123 = @('Raw_Spectral_Data_File','Raw Spectral Data File')
    124 = @('Derived_Spectral_Data_File','Derived Spectral Data File')
    125 = @('Peptide_Assignment_File','Peptide Assignment File')
    126 = @('Protein_Assignment_File','Protein Assignment File')
    127 = @('Post_trans_mode_ass_file','Post Translational Modification Assignment File')
    128 = @('NMR_Assay_Name','NMR Assay Name')
    129 = @('Free_Induct_Decay_Data_File_','Free Induction Decay Data File ')
    130 = @('Acq_Parameter_Data_File_','Acquisition Parameter Data File ')
    131 = @('Derived_Spectral_Data_File','Derived Spectral Data File')
}

# Insert a new column before column B ("entity"); this shifts the existing
# B..P columns to C..Q. Column A (the original "name" values) is untouched
# by the insert, so we still have the old text there for a moment.
$ws.Columns.Item(2).Insert()

# Header row: A1 stays "name", new B1 is "label"
$ws.Range("B1").Value = "label"

# For every data row: move the (still present) old A value into B, then
# overwrite A with the new machine-readable short name.
for ($r = 2; $r -le 131; $r++) {
    $pair = $rowData[$r]
    $oldName = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $oldName
    $ws.Cells.Item($r, 1).Value = $pair[0]
}

# Selection / active-sheet bookkeeping: "attributes" becomes the active tab
# with A2:A131 selected (was "packages" previously active).
$ws.Activate()
$ws.Range("A2:A131").Select()
